$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two columns (Active Residential, Log(Active Residential))
$ws.Range("D1:E8").Delete()

# Update the "Active and Vacant Residential" and "Log(...)" columns with new values
$ws.Range("B2").Value = "37.15 * (3.18)"
$ws.Range("C2").Value = "0.022 * (0.0021)"

$ws.Range("B3").Value = "52.73 * (5.83)"
$ws.Range("C3").Value = "0.0379 * (0.0045)"

$ws.Range("B4").Value = "60.03 * (11.05)"
$ws.Range("C4").Value = "0.0289  (0.0132)"

$ws.Range("B5").Value = "34.85 * (6.38)"
$ws.Range("C5").Value = "0.0228 * (0.0041)"

$ws.Range("B6").Value = "39.53 * (8.15)"
$ws.Range("C6").Value = "0.0181 * (0.0049)"

$ws.Range("B7").Value = "16.76 * (5.67)"
$ws.Range("C7").Value = "0.0117  (0.0057)"

$ws.Range("B8").Value = "19.61  (7.79)"
$ws.Range("C8").Value = "0.0066  (0.0045)"
